$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8328052759170532
$ws.Range("B1").Value = 1.363547682762146
$ws.Range("D1").Value = 1.749380826950073
$ws.Range("E1").Value = 1.143595337867737
